$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 3.781711156805759)
    3 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 6.740334628841572)
    4 = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 8.660232485948974, 1, 39.7764191927396)
    5 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    6 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    7 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 0, 11.80884686099532)
    8 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 0, 11.80884686099532)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
